$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the two "legajo" (student ID) numbers in the membership
#    paragraph. Each "00000" placeholder is disambiguated by its
#    surrounding name so the right one gets replaced.
# ------------------------------------------------------------------
$d.Content.Find.Execute("Martín Lucas, legajo: 00000", $true, $false, $false, $false, $false, $true, 1, $false, "Martín Lucas, legajo: 70409", 2) | Out-Null
$d.Content.Find.Execute("Sivoff Nicolas, legajo: 00000", $true, $false, $false, $false, $false, $true, 1, $false, "Sivoff Nicolas, legajo: 73841", 2) | Out-Null

# ------------------------------------------------------------------
# 2. Fix conjugation: "Se modificaron de las tablas" -> "Se modificó de las tablas"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Se modificaron de las tablas", $true, $false, $false, $false, $false, $true, 1, $false, "Se modificó de las tablas", 2) | Out-Null

# ------------------------------------------------------------------
# 3. Give the blank paragraphs surrounding the "El proyecto..." and
#    "Llevamos a cabo..." paragraphs an explicit es-ES run language
#    (matching the rest of the document's runs).
# ------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $ptext = $para.Range.Text
    if ($ptext -eq "`r" -and $para.Range.LanguageID -ne "es-ES") {
        $para.Range.LanguageID = "es-ES"
    }
}

# ------------------------------------------------------------------
# 4. Append two new bullet items to the closing list, after the
#    "Se modificó..." item, re-using the same ListParagraph style /
#    numbering that the rest of the list uses.
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("Se opto por la tecnología desarrollada por Microsoft " + [char]8220 + "OLE BD" + [char]8221 + " para acceder a la información contenida en la Base de Datos.")

$r2 = $d.Paragraphs.Last.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$d.Paragraphs.Last.Range.InsertAfter("El aspecto nombre, aspecto visual, y funcionamiento del sistema se realizo de manera consensuada por todos los integrantes.")
